$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New student row (row 4) - test case reproducing import when birth_date
# (column F) and option_ase (column X) are empty.
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = 11

$cText = $ws.Range("C3").Text
$ws.Range("C4").Value = $cText

$ws.Range("D4").Value = "Irène"
$ws.Range("E4").Value = "F"

# F4 keeps the date number format of F3 (birth date column) but is left
# empty - this reproduces the "birth_date empty" import scenario.
$ws.Range("F4").NumberFormat = $ws.Range("F3").NumberFormat

$gText = $ws.Range("G3").Text
$ws.Range("G4").Value = $gText

$hText = $ws.Range("H3").Text
$ws.Range("H4").Value = $hText

$iText = $ws.Range("I3").Text
$ws.Range("I4").Value = $iText
$ws.Range("I4").Font.Color = $ws.Range("I3").Font.Color

$kText = $ws.Range("K3").Text
$ws.Range("K4").Value = $kText

$nText = $ws.Range("N3").Text
$ws.Range("N4").Value = $nText

$oText = $ws.Range("O3").Text
$ws.Range("O4").Value = $oText

$uText = $ws.Range("U3").Text
$ws.Range("U4").Value = $uText

$vText = $ws.Range("V3").Text
$ws.Range("V4").Value = $vText

# option_ase (X4) intentionally left empty - reproduces "option_ase empty"
# import scenario (X3 on the prior row had a value).

$alText = $ws.Range("AL3").Text
$ws.Range("AL4").Value = $alText

# Move the active selection back to the top-left of the sheet / the new
# row, matching the refreshed view state.
$ws.Range("A1").Select() | Out-Null
$ws.Range("A4").Select() | Out-Null
